$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 and 43: coin metadata swapped (CEJI <-> BKEXToken) ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Price (column D) updates. Cells are stored as text, so force
#     the "Text" number format before assigning the numeric-looking
#     string, otherwise Excel would auto-convert it to a real number. ---
$priceUpdates = @{
    2 = "260.62"
    3 = "22.70"
    4 = "6.212"
    5 = "0.06087"
    6 = "3.510"
    7 = "6.710"
    8 = "1.360"
    9 = "0.7989"
    10 = "0.1575"
    11 = "0.08115"
    12 = "0.03311"
    13 = "0.03138"
    14 = "0.09264"
    15 = "3.890"
    16 = "0.001694"
    17 = "0.04830"
    18 = "0.0006204"
    19 = "0.006237"
    20 = "0.001100"
    21 = "0.003375"
    23 = "3.694"
    24 = "2.292"
    25 = "0.3375"
    27 = "0.0006169"
    41 = "0.007174"
    42 = "0.003902"
    43 = "0.1118"
    44 = "0.01020"
    46 = "0.00006037"
    47 = "0.00000000750"
    48 = "0.7004"
    49 = "0.05391"
    50 = "0.00001501"
}
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}
